# Updates odds/value columns for the week-of-2025-02-12 FlashScore sheet
# to reflect the latest odds snapshot (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("K2").Value = 2.38
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2
$ws.Range("W2").Value = 1.36
$ws.Range("X2").Value = 3
$ws.Range("Y2").Value = 1.95
$ws.Range("Z2").Value = 1.8
$ws.Range("AA2").Value = 19
$ws.Range("AD2").Value = 101
$ws.Range("AG2").Value = 11
$ws.Range("AK2").Value = 451
$ws.Range("AL2").Value = 7.5
# Row 3
$ws.Range("G3").Value = 2.2
$ws.Range("I3").Value = 3.7
$ws.Range("J3").Value = 3.1
$ws.Range("S3").Value = 4.8
$ws.Range("T3").Value = 1.19
$ws.Range("W3").Value = 1.62
$ws.Range("X3").Value = 2.2
$ws.Range("AG3").Value = 6
$ws.Range("AL3").Value = 7.5
$ws.Range("AM3").Value = 17
$ws.Range("AN3").Value = 15
$ws.Range("AR3").Value = 2.03
$ws.Range("AS3").Value = 1.83
# Row 4
$ws.Range("G4").Value = 2.75
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 2.75
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 2
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57
$ws.Range("T4").Value = 1.29
$ws.Range("W4").Value = 1.5
$ws.Range("X4").Value = 2.5
$ws.Range("AD4").Value = 29
$ws.Range("AK4").Value = 351
$ws.Range("AM4").Value = 12
$ws.Range("AP4").Value = 23
# Row 5
$ws.Range("G5").Value = 1.91
$ws.Range("I5").Value = 4.75
$ws.Range("J5").Value = 2.75
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 6
$ws.Range("AA5").Value = 4.75
$ws.Range("AL5").Value = 8.5
$ws.Range("AP5").Value = 51
$ws.Range("AQ5").Value = 67
$ws.Range("AR5").Value = 2.2
# Row 6
$ws.Range("G6").Value = 1.83
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 2.63
$ws.Range("K6").Value = 1.91
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6
$ws.Range("O6").Value = 1.53
$ws.Range("P6").Value = 2.38
$ws.Range("S6").Value = 4.6
$ws.Range("T6").Value = 1.2
$ws.Range("U6").Value = 5.5
$ws.Range("V6").Value = 1.14
$ws.Range("W6").Value = 1.57
$ws.Range("X6").Value = 2.25
$ws.Range("AB6").Value = 7
$ws.Range("AC6").Value = 9.5
$ws.Range("AD6").Value = 15
$ws.Range("AE6").Value = 19
$ws.Range("AI6").Value = 21
$ws.Range("AJ6").Value = 81
$ws.Range("AL6").Value = 9.5
$ws.Range("AM6").Value = 23
$ws.Range("AN6").Value = 17
$ws.Range("AP6").Value = 51
$ws.Range("AR6").Value = 2.03
$ws.Range("AS6").Value = 1.83
# Row 7
$ws.Range("G7").Value = 5.75
$ws.Range("H7").Value = 3.75
$ws.Range("I7").Value = 1.62
$ws.Range("J7").Value = 6.5
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 2.25
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.75
$ws.Range("Q7").Value = 2.25
$ws.Range("R7").Value = 1.62
$ws.Range("S7").Value = 3.45
$ws.Range("T7").Value = 1.32
$ws.Range("U7").Value = 4.33
$ws.Range("V7").Value = 1.2
$ws.Range("W7").Value = 1.5
$ws.Range("X7").Value = 2.5
$ws.Range("AA7").Value = 12
$ws.Range("AB7").Value = 29
$ws.Range("AC7").Value = 19
$ws.Range("AD7").Value = 67
$ws.Range("AN7").Value = 9
$ws.Range("AR7").Value = 1.7
$ws.Range("AS7").Value = 2.17
# Row 8
$ws.Range("G8").Value = 1.73
$ws.Range("I8").Value = 5.75
$ws.Range("J8").Value = 2.5
$ws.Range("M8").Value = 1.14
$ws.Range("N8").Value = 5.5
$ws.Range("O8").Value = 1.57
$ws.Range("P8").Value = 2.25
$ws.Range("Q8").Value = 2.88
$ws.Range("R8").Value = 1.4
$ws.Range("S8").Value = 4.5
$ws.Range("T8").Value = 1.19
$ws.Range("U8").Value = 6
$ws.Range("V8").Value = 1.13
$ws.Range("AA8").Value = 4.75
$ws.Range("AB8").Value = 6.5
$ws.Range("AC8").Value = 10
$ws.Range("AG8").Value = 5.5
$ws.Range("AN8").Value = 21
$ws.Range("AO8").Value = 67
$ws.Range("AR8").Value = 2.1
$ws.Range("AS8").Value = 1.78
# Row 9
$ws.Range("G9").Value = 1.4
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = 2
$ws.Range("L9").Value = 9.5
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 7
$ws.Range("O9").Value = 1.44
$ws.Range("P9").Value = 2.63
$ws.Range("Q9").Value = 2.35
$ws.Range("R9").Value = 1.57
$ws.Range("Y9").Value = 2.75
$ws.Range("Z9").Value = 1.4
$ws.Range("AA9").Value = 4.75
$ws.Range("AB9").Value = 5
$ws.Range("AC9").Value = 10
$ws.Range("AD9").Value = 8.5
$ws.Range("AG9").Value = 7
$ws.Range("AH9").Value = 8.5
$ws.Range("AI9").Value = 29
$ws.Range("AJ9").Value = 126
$ws.Range("AN9").Value = 29
$ws.Range("AO9").Value = 126
$ws.Range("AP9").Value = 81
$ws.Range("AQ9").Value = 101
$ws.Range("AR9").Value = 1.78
# Row 10
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 2.8
$ws.Range("I10").Value = 2.63
$ws.Range("J10").Value = 3.75
$ws.Range("L10").Value = 3.4
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 6
$ws.Range("O10").Value = 1.53
$ws.Range("P10").Value = 2.38
$ws.Range("Q10").Value = 2.7
$ws.Range("R10").Value = 1.44
$ws.Range("S10").Value = 4.3
$ws.Range("T10").Value = 1.21
$ws.Range("U10").Value = 5.5
$ws.Range("V10").Value = 1.14
$ws.Range("AB10").Value = 13
$ws.Range("AC10").Value = 12
$ws.Range("AE10").Value = 29
$ws.Range("AG10").Value = 6
$ws.Range("AH10").Value = 5.5
$ws.Range("AI10").Value = 17
$ws.Range("AM10").Value = 11
$ws.Range("AO10").Value = 26
$ws.Range("AR10").Value = 2.05
$ws.Range("AS10").Value = 1.8
# Row 15
$ws.Range("G15").Value = 1.75
$ws.Range("H15").Value = 3.5
$ws.Range("I15").Value = 4.75
$ws.Range("J15").Value = 2.4
$ws.Range("O15").Value = 1.33
$ws.Range("P15").Value = 3.4
$ws.Range("Q15").Value = 2.05
$ws.Range("R15").Value = 1.8
$ws.Range("AA15").Value = 7
$ws.Range("AF15").Value = 29
$ws.Range("AG15").Value = 9.5
$ws.Range("AL15").Value = 12
# Row 16
$ws.Range("AR16").Value = 1.83
$ws.Range("AS16").Value = 2.03
# Row 17
$ws.Range("G17").Value = 3.5
$ws.Range("I17").Value = 2.25
$ws.Range("J17").Value = 4.33
$ws.Range("L17").Value = 3.1
$ws.Range("M17").Value = 1.13
$ws.Range("N17").Value = 6
$ws.Range("O17").Value = 1.57
$ws.Range("P17").Value = 2.38
$ws.Range("Y17").Value = 2.25
$ws.Range("Z17").Value = 1.57
$ws.Range("AG17").Value = 6
$ws.Range("AO17").Value = 21
$ws.Range("AP17").Value = 23
$ws.Range("AR17").Value = 2.1
$ws.Range("AS17").Value = 1.78
# Row 18
$ws.Range("J18").Value = 2
$ws.Range("K18").Value = 2.25
$ws.Range("N18").Value = 9.5
$ws.Range("Q18").Value = 2.03
$ws.Range("R18").Value = 1.83
$ws.Range("U18").Value = 3.5
$ws.Range("V18").Value = 1.3
$ws.Range("AI18").Value = 21
# Row 19
$ws.Range("G19").Value = 2.63
$ws.Range("I19").Value = 2.6
$ws.Range("J19").Value = 3.4
$ws.Range("Q19").Value = 2
$ws.Range("R19").Value = 1.8
$ws.Range("AB19").Value = 13
$ws.Range("AM19").Value = 12
# Row 20
$ws.Range("Q20").Value = 1.98
$ws.Range("R20").Value = 1.88
$ws.Range("U20").Value = 3.4
$ws.Range("V20").Value = 1.33
# Row 21
$ws.Range("G21").Value = 1.9
$ws.Range("H21").Value = 3.4
$ws.Range("J21").Value = 2.6
$ws.Range("K21").Value = 2.1
$ws.Range("L21").Value = 4.75
$ws.Range("M21").Value = 1.07
$ws.Range("N21").Value = 9
$ws.Range("O21").Value = 1.36
$ws.Range("P21").Value = 3.2
$ws.Range("Q21").Value = 2.1
$ws.Range("R21").Value = 1.73
$ws.Range("U21").Value = 4
$ws.Range("V21").Value = 1.25
$ws.Range("W21").Value = 1.44
$ws.Range("X21").Value = 2.63
$ws.Range("Y21").Value = 1.91
$ws.Range("Z21").Value = 1.8
$ws.Range("AA21").Value = 6.5
$ws.Range("AC21").Value = 9
$ws.Range("AE21").Value = 17
$ws.Range("AF21").Value = 29
$ws.Range("AG21").Value = 8.5
$ws.Range("AH21").Value = 6.5
$ws.Range("AI21").Value = 17
$ws.Range("AK21").Value = 351
$ws.Range("AL21").Value = 11
$ws.Range("AN21").Value = 15
